$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 30
$ws.Range("G30").Value = 1.6
$ws.Range("H30").Value = 3.6
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 2.25
$ws.Range("M30").Value = 1.11
$ws.Range("N30").Value = 6.5
$ws.Range("O30").Value = 1.5
$ws.Range("R30").Value = 1.5
$ws.Range("T30").Value = 1.17
$ws.Range("U30").Value = 1.53
$ws.Range("V30").Value = 2.38
$ws.Range("AG30").Value = 26
$ws.Range("AJ30").Value = 12
$ws.Range("AK30").Value = 34
$ws.Range("AL30").Value = 23
$ws.Range("AN30").Value = 67
$ws.Range("AR30").Value = 1.9
$ws.Range("AS30").Value = 1.95

# Row 45
$ws.Range("G45").Value = 1.65
$ws.Range("I45").Value = 5.75
$ws.Range("J45").Value = 2.3
$ws.Range("L45").Value = 6
$ws.Range("M45").Value = 1.08
$ws.Range("N45").Value = 7.5
$ws.Range("O45").Value = 1.4
$ws.Range("P45").Value = 2.75
$ws.Range("Q45").Value = 2.25
$ws.Range("R45").Value = 1.62
$ws.Range("U45").Value = 1.5
$ws.Range("V45").Value = 2.5
$ws.Range("W45").Value = 2.2
$ws.Range("X45").Value = 1.62
$ws.Range("Y45").Value = 5.5
$ws.Range("AB45").Value = 12
$ws.Range("AE45").Value = 7.5
$ws.Range("AG45").Value = 21
$ws.Range("AH45").Value = 81
$ws.Range("AJ45").Value = 12
$ws.Range("AK45").Value = 26
$ws.Range("AL45").Value = 19
$ws.Range("AM45").Value = 67
$ws.Range("AN45").Value = 51
$ws.Range("AP45").Value = 3.35
$ws.Range("AQ45").Value = 1.32
$ws.Range("AR45").Value = 1.67
$ws.Range("AS45").Value = 2.16

# Row 114
$ws.Range("G114").Value = 1.75
$ws.Range("H114").Value = 4
$ws.Range("I114").Value = 4.2
$ws.Range("J114").Value = 2.38
$ws.Range("L114").Value = 5
$ws.Range("M114").Value = 1.05
$ws.Range("N114").Value = 11
$ws.Range("O114").Value = 1.3
$ws.Range("P114").Value = 3.4
$ws.Range("Q114").Value = 1.98
$ws.Range("R114").Value = 1.88
$ws.Range("W114").Value = 1.91
$ws.Range("X114").Value = 1.91
$ws.Range("Y114").Value = 7
$ws.Range("Z114").Value = 8
$ws.Range("AB114").Value = 13
$ws.Range("AC114").Value = 15
$ws.Range("AF114").Value = 7.5
$ws.Range("AG114").Value = 17
$ws.Range("AI114").Value = 301
$ws.Range("AJ114").Value = 11
$ws.Range("AK114").Value = 21
$ws.Range("AL114").Value = 13
$ws.Range("AM114").Value = 41
$ws.Range("AN114").Value = 34

# Row 115
$ws.Range("G115").Value = 1.67
$ws.Range("H115").Value = 4.1
$ws.Range("J115").Value = 2.25
$ws.Range("M115").Value = 1.02
$ws.Range("O115").Value = 1.17
$ws.Range("P115").Value = 4.33
$ws.Range("S115").Value = 2.62
$ws.Range("T115").Value = 1.41
$ws.Range("W115").Value = 1.67
$ws.Range("Y115").Value = 8
$ws.Range("AC115").Value = 13
$ws.Range("AD115").Value = 23
$ws.Range("AE115").Value = 13

# Row 117
$ws.Range("G117").Value = 1.57
$ws.Range("H117").Value = 3.65
$ws.Range("I117").Value = 5.6
$ws.Range("J117").Value = 2.12
$ws.Range("L117").Value = 5.5
$ws.Range("P117").Value = 3.2
$ws.Range("R117").Value = 1.85
$ws.Range("T117").Value = 1.34
$ws.Range("W117").Value = 1.78
$ws.Range("X117").Value = 1.83
$ws.Range("Y117").Value = 6.8
$ws.Range("Z117").Value = 7.4
$ws.Range("AA117").Value = 7.9
$ws.Range("AB117").Value = 11.5
$ws.Range("AC117").Value = 12.5
$ws.Range("AD117").Value = 25
$ws.Range("AE117").Value = 10.5
$ws.Range("AF117").Value = 7.3
$ws.Range("AG117").Value = 16
$ws.Range("AH117").Value = 70
$ws.Range("AI117").Value = 600
$ws.Range("AJ117").Value = 15.5
$ws.Range("AK117").Value = 35
$ws.Range("AL117").Value = 17.5
$ws.Range("AM117").Value = 120
$ws.Range("AN117").Value = 60
$ws.Range("AO117").Value = 55
